$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Mdk"
$ws.Range("C2").Value = "Tspan1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.098888
$ws.Range("H2").Value = 6.296664
$ws.Range("I2").Value = 0.1082453658858517
$ws.Range("J2").Value = 0.1082453658858517
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.476486
$ws.Range("N2").Value = 1.429458
$ws.Range("O2").Value = 0.1483751124393585
$ws.Range("P2").Value = 0.1483751124393586
$ws.Range("Q2").Value = 1.000090747568
$ws.Range("R2").Value = 9.000816728112001
$ws.Range("S2").Value = 0.01606091833435275
$ws.Range("T2").Value = 0.01606091833435274

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Mdk"
$ws.Range("C3").Value = "Tspan1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.098888
$ws.Range("H3").Value = 6.296664
$ws.Range("I3").Value = 0.1082453658858517
$ws.Range("J3").Value = 0.1082453658858517
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.853664666666667
$ws.Range("N3").Value = 5.560994
$ws.Range("O3").Value = 0.5772209536933566
$ws.Range("P3").Value = 0.5772209536933566
$ws.Range("Q3").Value = 3.890634524890667
$ws.Range("R3").Value = 35.015710724016
$ws.Range("S3").Value = 0.06248149332951763
$ws.Range("T3").Value = 0.06248149332951761

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Mdk"
$ws.Range("C4").Value = "Tspan1"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.098888
$ws.Range("H4").Value = 6.296664
$ws.Range("I4").Value = 0.1082453658858517
$ws.Range("J4").Value = 0.1082453658858517
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.8812099999999999
$ws.Range("N4").Value = 2.64363
$ws.Range("O4").Value = 0.2744039338672849
$ws.Range("P4").Value = 0.2744039338672849
$ws.Range("Q4").Value = 1.84956109448
$ws.Range("R4").Value = 16.64604985032
$ws.Range("S4").Value = 0.0297029542219813
$ws.Range("T4").Value = 0.02970295422198129

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Mdk"
$ws.Range("C5").Value = "Tspan1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 15.87514366666667
$ws.Range("H5").Value = 47.625431
$ws.Range("I5").Value = 0.8187243600843848
$ws.Range("J5").Value = 0.8187243600843847
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.476486
$ws.Range("N5").Value = 1.429458
$ws.Range("O5").Value = 0.1483751124393585
$ws.Range("P5").Value = 0.1483751124393586
$ws.Range("Q5").Value = 7.564283705155334
$ws.Range("R5").Value = 68.07855334639801
$ws.Range("S5").Value = 0.1214783189843624
$ws.Range("T5").Value = 0.1214783189843625

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Mdk"
$ws.Range("C6").Value = "Tspan1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 15.87514366666667
$ws.Range("H6").Value = 47.625431
$ws.Range("I6").Value = 0.8187243600843848
$ws.Range("J6").Value = 0.8187243600843847
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.853664666666667
$ws.Range("N6").Value = 5.560994
$ws.Range("O6").Value = 0.5772209536933566
$ws.Range("P6").Value = 0.5772209536933566
$ws.Range("Q6").Value = 29.42719289315711
$ws.Range("R6").Value = 264.844736038414
$ws.Range("S6").Value = 0.4725848559398916
$ws.Range("T6").Value = 0.4725848559398916

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Mdk"
$ws.Range("C7").Value = "Tspan1"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 15.87514366666667
$ws.Range("H7").Value = 47.625431
$ws.Range("I7").Value = 0.8187243600843848
$ws.Range("J7").Value = 0.8187243600843847
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.8812099999999999
$ws.Range("N7").Value = 2.64363
$ws.Range("O7").Value = 0.2744039338672849
$ws.Range("P7").Value = 0.2744039338672849
$ws.Range("Q7").Value = 13.98933535050333
$ws.Range("R7").Value = 125.90401815453
$ws.Range("S7").Value = 0.2246611851601307
$ws.Range("T7").Value = 0.2246611851601306

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Mdk"
$ws.Range("C8").Value = "Tspan1"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.416064
$ws.Range("H8").Value = 4.248192
$ws.Range("I8").Value = 0.07303027402976368
$ws.Range("J8").Value = 0.07303027402976367
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.476486
$ws.Range("N8").Value = 1.429458
$ws.Range("O8").Value = 0.1483751124393585
$ws.Range("P8").Value = 0.1483751124393586
$ws.Range("Q8").Value = 0.6747346711039999
$ws.Range("R8").Value = 6.072612039936
$ws.Range("S8").Value = 0.01083587512064335
$ws.Range("T8").Value = 0.01083587512064335

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Mdk"
$ws.Range("C9").Value = "Tspan1"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.416064
$ws.Range("H9").Value = 4.248192
$ws.Range("I9").Value = 0.07303027402976368
$ws.Range("J9").Value = 0.07303027402976367
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1.853664666666667
$ws.Range("N9").Value = 5.560994
$ws.Range("O9").Value = 0.5772209536933566
$ws.Range("P9").Value = 0.5772209536933566
$ws.Range("Q9").Value = 2.624907802538667
$ws.Range("R9").Value = 23.624170222848
$ws.Range("S9").Value = 0.04215460442394736
$ws.Range("T9").Value = 0.04215460442394735

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Mdk"
$ws.Range("C10").Value = "Tspan1"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.416064
$ws.Range("H10").Value = 4.248192
$ws.Range("I10").Value = 0.07303027402976368
$ws.Range("J10").Value = 0.07303027402976367
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.8812099999999999
$ws.Range("N10").Value = 2.64363
$ws.Range("O10").Value = 0.2744039338672849
$ws.Range("P10").Value = 0.2744039338672849
$ws.Range("Q10").Value = 1.24784975744
$ws.Range("R10").Value = 11.23064781696
$ws.Range("S10").Value = 0.02003979448517297
$ws.Range("T10").Value = 0.02003979448517296
